# Weekly update: insert a new price record for "Vega Modelo de Temuco - Jengibre"
# right above the current row 322, pushing all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 322 (rows 322:341 shift down to 323:342).
$ws.Rows("322:322").Insert()

# Populate the newly inserted row with this week's data point.
$ws.Range("A322").Value = 10
$ws.Range("B322").Value = "Vega Modelo de Temuco"
$ws.Range("C322").Value = "La Araucanía"
$ws.Range("D322").Value = 45147
$ws.Range("E322").Value = 9
$ws.Range("F322").Value = 100114007
$ws.Range("G322").Value = "Jengibre"
$ws.Range("H322").Value = "Sin especificar"
$ws.Range("I322").Value = "Primera"
$ws.Range("J322").Value = 25
$ws.Range("K322").Value = 24000
$ws.Range("L322").Value = 24000
$ws.Range("M322").Value = 24000
$ws.Range("N322").Value = "$/caja 13 kilos"
$ws.Range("O322").Value = "Perú"
$ws.Range("P322").Value = 1846
$ws.Range("Q322").Value = 13
$ws.Range("R322").Value = "Hortaliza"
